# Normalizar ARTICULOS / Cierres
$d = $word.ActiveDocument

# 1) Mark the document "header" paragraphs (date line, ORDENANZA title,
#    VISTO:, CONSIDERANDO:) as outline level 1 (w:outlineLvl val="0").
$d.Paragraphs(1).OutlineLevel = 1
$d.Paragraphs(2).OutlineLevel = 1
$d.Paragraphs(3).OutlineLevel = 1
$d.Paragraphs(5).OutlineLevel = 1

# 2) Drop the stray leading-space run that precedes the body text of the
#    "VISTO:" paragraph.
$p4 = $d.Paragraphs(4)
$start = $p4.Range.Start
$lead = $d.Range($start, $start + 1)
if ($lead.Text -eq " ") {
    $lead.Delete()
}

# 3) Drop the stray leading-space run that precedes the body text of the
#    "CONSIDERANDO:" paragraph.
$p6 = $d.Paragraphs(6)
$start = $p6.Range.Start
$lead = $d.Range($start, $start + 1)
if ($lead.Text -eq " ") {
    $lead.Delete()
}

# 4) Trim the closing paragraph (ARTICULO OCTAVO) so it ends right after
#    "CÓPIESE y ARCHÍVESE." — drop the trailing clause about the Ley Nº 5529.
$d.Content.Find.Execute(
    ". Y en un todo de acuerdo a las facultades conferidas a los H. Concejos Deliberantes por la Ley Nº 5529; ",
    $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)
